# fix shark double counts
# Update landings/status figures on the area-specific sheets to correct
# a double-counting bug.

$wb = $excel.ActiveWorkbook

# --- "Status by Landings (Area)" sheet ---
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")
$wsArea.Range("C3").Value = 0.5344385436970592
$wsArea.Range("C5").Value = 0.6201980835478056
$wsArea.Range("C7").Value = 6.341402055340052
$wsArea.Range("C8").Value = 39.51851520369351
$wsArea.Range("C9").Value = 54.14008274096643
$wsArea.Range("C10").Value = 45.85991725903357
$wsArea.Range("C11").Value = 54.14008274096643

# --- "Status by Landings (Tier)" sheet ---
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")
$wsTier.Range("C4").Value = 0.5344385436970592
$wsTier.Range("E4").Value = 0.6201980835478056
$wsTier.Range("G4").Value = 6.341402055340052
$wsTier.Range("H4").Value = 39.51851520369351
$wsTier.Range("I4").Value = 54.14008274096643
$wsTier.Range("J4").Value = 45.85991725903357
$wsTier.Range("K4").Value = 54.14008274096643

$wsTier.Range("C5").Value = 0.5344385436970592
$wsTier.Range("E5").Value = 0.6201980835478056
$wsTier.Range("G5").Value = 6.341402055340052
$wsTier.Range("H5").Value = 39.51851520369351
$wsTier.Range("I5").Value = 54.14008274096643
$wsTier.Range("J5").Value = 45.85991725903357
$wsTier.Range("K5").Value = 54.14008274096643

# --- "Comparison by Landings" sheet ---
$wsComp = $wb.Worksheets.Item("Comparison by Landings")
$wsComp.Range("C2").Value = 95.38283744105355
$wsComp.Range("C3").Value = 6.341402055340052
$wsComp.Range("C4").Value = 39.51851520369351
$wsComp.Range("C5").Value = 54.14008274096643
$wsComp.Range("C6").Value = 45.85991725903357
$wsComp.Range("C7").Value = 54.14008274096643
